$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.021814
$ws.Range("H2").Value = 0.065442
$ws.Range("I2").Value = 0.1008129179549036
$ws.Range("J2").Value = 0.1008129179549036
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.006612826844
$ws.Range("R2").Value = 0.05951544159599999
$ws.Range("S2").Value = 0.000518466240712741
$ws.Range("T2").Value = 0.0005184662407127411

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.021814
$ws.Range("H3").Value = 0.065442
$ws.Range("I3").Value = 0.1008129179549036
$ws.Range("J3").Value = 0.1008129179549036
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("Q3").Value = 0.3700335360366667
$ws.Range("R3").Value = 3.33030182433
$ws.Range("S3").Value = 0.02901178284150051
$ws.Range("T3").Value = 0.02901178284150051

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.021814
$ws.Range("H4").Value = 0.065442
$ws.Range("I4").Value = 0.1008129179549036
$ws.Range("J4").Value = 0.1008129179549036
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 0.9091815613399999
$ws.Range("R4").Value = 8.182634052059999
$ws.Range("S4").Value = 0.07128266887269037
$ws.Range("T4").Value = 0.07128266887269037

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.194567
$ws.Range("H5").Value = 0.583701
$ws.Range("I5").Value = 0.8991870820450963
$ws.Range("J5").Value = 0.8991870820450963
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 0.058982207782
$ws.Range("R5").Value = 0.530839870038
$ws.Range("S5").Value = 0.0046243889729878
$ws.Range("T5").Value = 0.004624388972987801

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.194567
$ws.Range("H6").Value = 0.583701
$ws.Range("I6").Value = 0.8991870820450963
$ws.Range("J6").Value = 0.8991870820450963
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("Q6").Value = 3.300463693318334
$ws.Range("R6").Value = 29.704173239865
$ws.Range("S6").Value = 0.2587666430788589
$ws.Range("T6").Value = 0.2587666430788589

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.194567
$ws.Range("H7").Value = 0.583701
$ws.Range("I7").Value = 0.8991870820450963
$ws.Range("J7").Value = 0.8991870820450963
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 8.109321025270001
$ws.Range("R7").Value = 72.98388922743
$ws.Range("S7").Value = 0.6357960499932497
$ws.Range("T7").Value = 0.6357960499932497
